# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure D-column price cells keep their original text formatting
# (values like '526.84' must stay as text, matching the source inlineStr cells)
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '60.659.41'
$ws.Range('E2').Value = '  -1.68%  '
$ws.Range('D3').Value = '2.903.16'
$ws.Range('E3').Value = '  -2.60%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '526.84'
$ws.Range('E5').Value = '  -2.64%  '
$ws.Range('D6').Value = '143.43'
$ws.Range('E6').Value = '  -5.52%  '
$ws.Range('D8').Value = '0.548'
$ws.Range('E8').Value = '  -3.36%  '
$ws.Range('D9').Value = '2.907.72'
$ws.Range('E9').Value = '  -2.82%  '
$ws.Range('E10').Value = '  -4.93%  '
$ws.Range('D11').Value = '5.99'
$ws.Range('E11').Value = '  -2.65%  '
$ws.Range('E12').Value = '  -2.84%  '
$ws.Range('D13').Value = '3.408.61'
$ws.Range('E13').Value = '  -2.73%  '
$ws.Range('D14').Value = '0.128'
$ws.Range('E14').Value = '  +3.19%  '
$ws.Range('D15').Value = '60.622.62'
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('D16').Value = '22.58'
$ws.Range('E16').Value = '  -5.58%  '
$ws.Range('D17').Value = '2.907.13'
$ws.Range('E17').Value = '  -2.51%  '
$ws.Range('E18').Value = '  -4.30%  '
$ws.Range('E19').Value = '  -3.91%  '
$ws.Range('D20').Value = '11.57'
$ws.Range('E20').Value = '  -4.00%  '
$ws.Range('D21').Value = '350.77'
$ws.Range('E21').Value = '  -8.18%  '
$ws.Range('D22').Value = '6.54'
$ws.Range('E22').Value = '  -2.61%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').Value = '5.72'
$ws.Range('E24').Value = '  +1.31%  '
$ws.Range('D25').Value = '64.82'
$ws.Range('E25').Value = '  -1.67%  '
$ws.Range('D26').Value = '0.451'
$ws.Range('E26').Value = '  -4.40%  '
$ws.Range('E27').Value = '  -6.91%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('D29').Value = '7.83'
$ws.Range('E29').Value = '  -4.30%  '
$ws.Range('D30').Value = '0.0₃0857'
$ws.Range('E30').Value = '  -8.89%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  -2.92%  '
$ws.Range('D33').Value = '19.53'
$ws.Range('E33').Value = '  -4.82%  '
$ws.Range('D34').Value = '152.17'
$ws.Range('E34').Value = '  -5.11%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '4.30'
$ws.Range('E35').Value = '  -6.10%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').Value = '5.56'
$ws.Range('E36').Value = '  -5.95%  '
$ws.Range('D37').Value = '0.995'
$ws.Range('E37').Value = '  -7.00%  '
$ws.Range('E38').Value = '  -5.90%  '
$ws.Range('D39').Value = '37.57'
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').Value = '1.47'
$ws.Range('E40').Value = '  -5.13%  '
$ws.Range('D41').Value = '3.71'
$ws.Range('E41').Value = '  -4.86%  '
$ws.Range('D42').Value = '2.289.87'
$ws.Range('E42').Value = '  -5.49%  '
$ws.Range('E43').Value = '  -3.44%  '
$ws.Range('D44').Value = '0.0579'
$ws.Range('E44').Value = '  -1.88%  '
$ws.Range('D45').Value = '20.39'
$ws.Range('E45').Value = '  -7.47%  '
$ws.Range('D46').Value = '0.997'
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').Value = '4.93'
$ws.Range('E47').Value = '  -4.58%  '
$ws.Range('D48').Value = '0.0237'
$ws.Range('E48').Value = '  -3.20%  '
$ws.Range('D50').Value = '0.0916'
$ws.Range('E50').Value = '  -3.97%  '
$ws.Range('D51').Value = '18.36'
$ws.Range('E51').Value = '  -7.21%  '
